$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 207, shifting rows 207:298 down to 208:299.
$ws.Rows.Item(207).Insert()

# Populate the newly inserted row 207 with data (copy of row structure with new values).
$ws.Cells.Item(207, 1).Value = 4
$ws.Cells.Item(207, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(207, 3).Value = "Los Lagos"
$ws.Cells.Item(207, 4).Value = 44845
$ws.Cells.Item(207, 5).Value = 10
$ws.Cells.Item(207, 6).Value = 100112044
$ws.Cells.Item(207, 7).Value = "Perejil"
$ws.Cells.Item(207, 8).Value = "Sin especificar"
$ws.Cells.Item(207, 9).Value = "Primera"
$ws.Cells.Item(207, 10).Value = 160
$ws.Cells.Item(207, 11).Value = 5000
$ws.Cells.Item(207, 12).Value = 5000
$ws.Cells.Item(207, 13).Value = 5000
$ws.Cells.Item(207, 14).Value = "$/docena de atados (3 kilos)"
$ws.Cells.Item(207, 15).Value = "Región Metropolitana"
$ws.Cells.Item(207, 16).Value = 1667
$ws.Cells.Item(207, 17).Value = 3
$ws.Cells.Item(207, 18).Value = "Hortaliza"
